$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @("2023-12-07 10:21:33", 0.0006000000000000001),
    @("2023-12-07 10:21:42", 0.0004),
    @("2023-12-07 10:21:50", 0.0006000000000000001),
    @("2023-12-07 10:22:02", 0.0004)
)

$startRow = 55
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $newRows[$i][0]
    $ws.Cells.Item($r, 2).Value = $newRows[$i][1]
}
